# Update Name of Algo - adjust a handful of imputed values in the
# RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -20.40989999999997
$ws.Range("A12").Value = -22.38260000000004
$ws.Range("C13").Value = -12.77679999999999
$ws.Range("A18").Value = -22.31730000000003
